$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row (row 11): correct marks value 3 -> 5
$ws.Range("B11").Value = 5

# Update the "Total" row (row 12): total marks value 48 -> 80
$ws.Range("B12").Value = 80

# Update corresponding displayed fraction "40/84" -> "80/140"
$ws.Range("E12").Value = "80/140"
